$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 32 ("DRAIAMSSO003/LogIn" row),
# shifting the existing rows 32-35 down to 34-37.
$ws.Rows("32:33").Insert()

# Copy the formatting (borders / wrap-text) of row 31 (the row immediately
# above the insertion point) onto the two freshly-inserted rows so they
# match the rest of the table's look (thin box border all round, column C
# wrapped).
$srcRow = $ws.Range("A31:E31")
$newRows = $ws.Range("A32:E33")
$newRows.Borders.LineStyle = 1
$newRows.Borders.Color = 0
$ws.Range("C32:C33").WrapText = $true

# Row 32: DRAIAM112 / OP11 / Suspend User / Y
$ws.Range("A32").Value = "DRAIAM112"
$ws.Range("B32").Value = "OP11"
$ws.Range("C32").Value = "Suspend User"
$ws.Range("D32").Value = "Y"

# Row 33: DRAIAM113 / OP113 / Existing User / Y
$ws.Range("A33").Value = "DRAIAM113"
$ws.Range("B33").Value = "OP113"
$ws.Range("C33").Value = "Existing User"
$ws.Range("D33").Value = "Y"
